$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 56 updates ---
# E56: Transaction value  " 0.185USDT" -> "         0.185  USDT"
$ws.Range("E56").Value = "         0.185  USDT"
# H56: Status  "IN PROGRESS" -> "CANCEL"
$ws.Range("H56").Value = "CANCEL"
# I56: Finalized date  " " -> "2017-05-13-20:35:22 "
$ws.Range("I56").Value = "2017-05-13-20:35:22 "

# --- Row 57 (new row) ---
# Copy date-style formatting (style index 2) onto A57 / I57 first
$ws.Range("A2").Copy()
$ws.Range("A57").PasteSpecial(-4122)
$ws.Range("A57").Value = 42871.274247685185

$ws.Range("B57").Value = "            Buy"
$c = $ws.Range("B57").Characters(13, 3)
$c.Font.Color = 5287936

$ws.Range("C57").Value = "        LTC"
$ws.Range("D57").Value = 24.043316999999998
$ws.Range("E57").Value = "            23.8 USDT"
$ws.Range("F57").Value = "        1.390 LTC"
$ws.Range("G57").Value = " LTC/USDT0000002"
$ws.Range("H57").Value = "IN PROGRESS"

$ws.Range("A2").Copy()
$ws.Range("I57").PasteSpecial(-4122)

$ws.Range("K57").Value = "   "

# --- view/selection state ---
$ws.Range("B62").Select()
